# Update column G ("K" - strikeouts) with regenerated values for each row.
# Per commit message: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" -- the only cell-level data change
# visible in the diff is the per-row "K" (strikeout) count in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 3
    6 = 1
    7 = 2
    8 = 3
    9 = 1
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 3
    19 = 3
    21 = 3
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 3
    34 = 0
    35 = 3
    36 = 0
    37 = 2
    38 = 0
    39 = 1
    40 = 4
    41 = 1
    42 = 2
    43 = 1
    44 = 4
    45 = 2
    46 = 3
    47 = 1
    48 = 3
    49 = 2
    50 = 4
    51 = 3
    52 = 1
    54 = 1
    55 = 5
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
